$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Plan" text for weeks 4, 5 and 6 (column C) ---
# Shared-string table order follows the order the values were typed in:
# "Review presentation 1", "Review presentation 2", then "Draft presentation".
$ws.Range("C6").Value = "Review presentation 1"
$ws.Range("C7").Value = "Review presentation 2"
$ws.Range("C5").Value = "Draft presentation"

# --- Match the grey header fill already used on the rest of the row ---
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- Highlight the PRESENTATION day (G7) with an accent colour fill ---
# Pre-seed a solid pattern fill (re-using the existing grey fill) so the
# following ThemeColor assignment mutates the fill in place instead of
# allocating a spurious intermediate fill entry.
$ws.Range("E2").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Interior.ThemeColor = 8

# --- Give the FINAL SUBMISSION cell (E14) the same accent colour fill,
#     keeping it border-free like it was before ---
$ws.Range("E2").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Interior.ThemeColor = 8
$ws.Range("E14").Borders.LineStyle = -4142

$excel.CutCopyMode = 0

# --- Scroll / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E17").Select()
